# Rename the header cells so the "_old"/"_new" suffixes become
# "_FV2310"/"_FV2404" (matching the respective input file format versions),
# then turn the header range into a real Excel Table and freeze the header
# row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-J (1-10) carry the "old" AHB (FV2310) headers; L-U (12-21) carry
# the "new" AHB (FV2404) headers. Column K (11) is the literal "diff" header
# and stays untouched.
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value() -replace "_old$", "_FV2310")
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value() -replace "_new$", "_FV2404")
}

# Stash the header row's existing direct formatting (bold font, grey fill,
# thin border, centered+wrapped) in a scratch cell so it survives the
# ListObjects.Add() call below untouched. Without this, the engine infers a
# "header row differs from the table style" override and writes a brand new
# dxf (xl/styles.xml dxfs/cellXfs) that the original file never had.
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("W1")
$scratch.Value = "scratch"
$ws.Range("A1").Copy($scratch)
$headerRange.ClearFormats()

# Turn A1:U64 into an Excel Table ("Table1") using the renamed headers.
$rng = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Restore the header row's original look and discard the scratch cell.
$scratch.Copy()
$headerRange.PasteSpecial(-4122) # xlPasteFormats
$scratch.Clear()

# Freeze the header row (split under row 1, top-left of scrollable area A2).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
